$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Columns("H").Insert()
$ws.Columns("K").Insert()

$u = $ws.Range("N10:Y12")
$u = $excel.Union($u, $ws.Range("O9:T9"))
$u = $excel.Union($u, $ws.Range("N3:X3"))
$u = $excel.Union($u, $ws.Range("N8:T8"))
$u = $excel.Union($u, $ws.Range("N7:U7"))
$u = $excel.Union($u, $ws.Range("N6:V6"))
$u = $excel.Union($u, $ws.Range("N4:W5"))
Write-Host "Union Address=$($u.Address())"
Write-Host "Areas.Count=$($u.Areas().Count())"

$cf2 = $ws.Cells.FormatConditions
Write-Host "Count=$($cf2.Count())"
for ($i = 1; $i -le $cf2.Count(); $i++) {
    $c = $cf2.Item($i)
    Write-Host "Item $i Type=$($c.Type()) Priority=$($c.Priority())"
    $a = $c.AppliesTo()
    Write-Host "AppliesTo=$($a.Address())"
}

for ($i = 5; $i -le 6; $i++) {
    $c = $cf2.Item($i)
    $c.ModifyAppliesToRange($u)
}

Write-Host "--- after modify ---"
for ($i = 1; $i -le $cf2.Count(); $i++) {
    $c = $cf2.Item($i)
    $a = $c.AppliesTo()
    Write-Host "Item $i AppliesTo=$($a.Address())"
}
